$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "001"
$ws.Range("C2").Value = "001"
$ws.Range("B4").Value = "001"
$ws.Range("C4").Value = "002"

$ws.Range("O15").Select()
